$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.186828970909119
$ws.Range("B1").Value = 2.117293834686279
$ws.Range("C1").Value = 4.309908390045166
$ws.Range("D1").Value = 2.942328929901123
$ws.Range("E1").Value = 1.216072201728821
